$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '55.815.54'
$ws.Cells.Item(2, 5).Value = '  +9.20%  '
$ws.Cells.Item(3, 4).Value = '2.518.87'
$ws.Cells.Item(3, 5).Value = '  +13.30%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '491.08'
$ws.Cells.Item(5, 5).Value = '  +16.93%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '141.93'
$ws.Cells.Item(6, 5).Value = '  +23.36%  '
$ws.Cells.Item(7, 5).Value = '  +0.40%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.516'
$ws.Cells.Item(8, 5).Value = '  +14.49%  '
$ws.Cells.Item(9, 4).Value = '2.513.79'
$ws.Cells.Item(9, 5).Value = '  +13.46%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0994'
$ws.Cells.Item(10, 5).Value = '  +16.47%  '
$ws.Cells.Item(11, 5).Value = '  +10.85%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.331'
$ws.Cells.Item(12, 5).Value = '  +15.31%  '
$ws.Cells.Item(13, 5).Value = '  +3.62%  '
$ws.Cells.Item(14, 4).Value = '2.955.52'
$ws.Cells.Item(14, 5).Value = '  +15.64%  '
$ws.Cells.Item(15, 4).Value = '55.853.64'
$ws.Cells.Item(15, 5).Value = '  +8.82%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '20.88'
$ws.Cells.Item(16, 5).Value = '  +15.30%  '
$ws.Cells.Item(17, 5).Value = '  +23.79%  '
$ws.Cells.Item(18, 4).Value = '2.520.79'
$ws.Cells.Item(18, 5).Value = '  +14.11%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.43'
$ws.Cells.Item(19, 5).Value = '  +17.71%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '324.06'
$ws.Cells.Item(20, 5).Value = '  +13.17%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.08'
$ws.Cells.Item(21, 5).Value = '  +18.91%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '1.00'
$ws.Cells.Item(22, 5).Value = '  +0.17%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.79'
$ws.Cells.Item(23, 5).Value = '  +14.72%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '58.55'
$ws.Cells.Item(24, 5).Value = '  +12.82%  '
$ws.Cells.Item(25, 5).Value = '  +26.19%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.412'
$ws.Cells.Item(26, 5).Value = '  +18.53%  '
$ws.Cells.Item(27, 5).Value = '  +0.40%  '
$ws.Cells.Item(28, 4).Value = '2.631.28'
$ws.Cells.Item(28, 5).Value = '  +12.67%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.53'
$ws.Cells.Item(29, 5).Value = '  +14.49%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0808'
$ws.Cells.Item(30, 5).Value = '  +25.46%  '
$ws.Cells.Item(31, 5).Value = '  +0.46%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '150.21'
$ws.Cells.Item(32, 5).Value = '  +7.49%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '18.35'
$ws.Cells.Item(33, 5).Value = '  +11.88%  '
$ws.Cells.Item(34, 5).Value = '  +19.43%  '
$ws.Cells.Item(35, 5).Value = '  +16.93%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.896'
$ws.Cells.Item(36, 5).Value = '  +14.59%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '3.75'
$ws.Cells.Item(37, 5).Value = '  +13.65%  '
$ws.Cells.Item(38, 5).Value = '  +18.50%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '34.34'
$ws.Cells.Item(39, 5).Value = '  +10.86%  '
$ws.Cells.Item(40, 5).Value = '  +24.49%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0560'
$ws.Cells.Item(41, 5).Value = '  +18.19%  '
$ws.Cells.Item(42, 5).Value = '  +0.45%  '
$ws.Cells.Item(43, 5).Value = '  +15.69%  '
$ws.Cells.Item(44, 5).Value = '  +15.20%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '4.77'
$ws.Cells.Item(45, 5).Value = '  +30.64%  '
$ws.Cells.Item(46, 4).Value = '2.010.28'
$ws.Cells.Item(46, 5).Value = '  +10.69%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '258.92'
$ws.Cells.Item(47, 5).Value = '  +49.17%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0915'
$ws.Cells.Item(48, 5).Value = '  +15.09%  '
$ws.Cells.Item(49, 2).Value = 'VeChain'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0227'
$ws.Cells.Item(49, 5).Value = '  +15.73%  '
$ws.Cells.Item(50, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '10.12'
$ws.Cells.Item(50, 5).Value = '  -0.53%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '17.78'
$ws.Cells.Item(51, 5).Value = '  +18.98%  '
